# Scheduled price-refresh run: update computed profit columns (H..N) on the
# per-job-class Leve tables after a new Market Board price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 40589.05
$ws.Range("I17").Value = 200000
$ws.Range("J17").Value = 37840.586
$ws.Range("K17").Value = 600000
$ws.Range("L17").Value = 113521.758
$ws.Range("M17").Value = -599832
$ws.Range("N17").Value = -113857.758

# Row 26: Everything Is Impossible / Budding Ash Wand
$ws.Range("H26").Value = 48338.332
$ws.Range("J26").Value = 48338.332
$ws.Range("L26").Value = 48338.332
$ws.Range("N26").Value = -49026.332

# Row 64: Forged from the Void / Void Glue
$ws.Range("H64").Value = 3985
$ws.Range("I64").Value = 3988
$ws.Range("J64").Value = 3979
$ws.Range("K64").Value = 3988
$ws.Range("L64").Value = 3979
$ws.Range("M64").Value = -3740
$ws.Range("N64").Value = -4475

# Row 67: Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 3985
$ws.Range("I67").Value = 3988
$ws.Range("J67").Value = 3979
$ws.Range("K67").Value = 3988
$ws.Range("L67").Value = 3979
$ws.Range("M67").Value = -3130
$ws.Range("N67").Value = -5695

# Row 111: An Eye for Healing / Grade 1 Dexterity Alkahest
$ws.Range("H111").Value = 2249.5
$ws.Range("I111").Value = 2249.5
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 6748.5
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -3681.5
$ws.Range("N111").ClearContents()

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 1284.898
$ws.Range("I112").Value = 400
$ws.Range("J112").Value = 1342.6086
$ws.Range("K112").Value = 1200
$ws.Range("L112").Value = 4027.8258
$ws.Range("M112").Value = -92
$ws.Range("N112").Value = -6243.825800000001

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2956.9429
$ws.Range("I137").Value = 2966.4333
$ws.Range("J137").Value = 2900
$ws.Range("K137").Value = 8899.2999
$ws.Range("L137").Value = 8700
$ws.Range("M137").Value = -6349.2999
$ws.Range("N137").Value = -13800

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 3847.91
$ws.Range("I32").Value = 3297.5889
$ws.Range("J32").Value = 8800.799999999999
$ws.Range("K32").Value = 3297.5889
$ws.Range("L32").Value = 8800.799999999999
$ws.Range("M32").Value = -3010.5889
$ws.Range("N32").Value = -9374.799999999999

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 3464.9167
$ws.Range("I122").Value = 2259.875
$ws.Range("J122").Value = 5875
$ws.Range("K122").Value = 6779.625
$ws.Range("L122").Value = 17625
$ws.Range("M122").Value = -4329.625
$ws.Range("N122").Value = -22525

# Row 137: Odd Instruments / Cobalt Tungsten Alembic
$ws.Range("H137").Value = 29642.857
$ws.Range("J137").Value = 29642.857
$ws.Range("L137").Value = 29642.857
$ws.Range("N137").Value = -39842.857

# Row 139: Backing up My Words / Titanium Gold Thornplate of Fending
$ws.Range("H139").Value = 29546.316
$ws.Range("J139").Value = 29546.316
$ws.Range("L139").Value = 29546.316
$ws.Range("N139").Value = -39826.316

$ws = $wb.Worksheets.Item("CRP")
# Row 20: Re-crating the Scene / Iron Spear
$ws.Range("H20").Value = 44999.5
$ws.Range("J20").Value = 44999.5
$ws.Range("L20").Value = 44999.5
$ws.Range("N20").Value = -45471.5

# Row 30: Polearms Aplenty / Iron Spear
$ws.Range("H30").Value = 44999.5
$ws.Range("J30").Value = 44999.5
$ws.Range("L30").Value = 44999.5
$ws.Range("N30").Value = -45181.5

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3788.5117
$ws.Range("I31").Value = 2908.5
$ws.Range("J31").Value = 4900.1055
$ws.Range("K31").Value = 2908.5
$ws.Range("L31").Value = 4900.1055
$ws.Range("M31").Value = -2613.5
$ws.Range("N31").Value = -5490.1055

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3788.5117
$ws.Range("I34").Value = 2908.5
$ws.Range("J34").Value = 4900.1055
$ws.Range("K34").Value = 2908.5
$ws.Range("L34").Value = 4900.1055
$ws.Range("M34").Value = -2706.5
$ws.Range("N34").Value = -5304.1055

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 3366.2
$ws.Range("I122").Value = 3028.25
$ws.Range("J122").Value = 3752.4285
$ws.Range("K122").Value = 9084.75
$ws.Range("L122").Value = 11257.2855
$ws.Range("M122").Value = -6634.75
$ws.Range("N122").Value = -16157.2855

# Row 128: An A-prop-riate Request / Ironwood Spear
$ws.Range("H128").Value = 44999.5
$ws.Range("J128").Value = 44999.5
$ws.Range("L128").Value = 44999.5
$ws.Range("N128").Value = -54959.5

$ws = $wb.Worksheets.Item("CUL")
# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 2112.1667
$ws.Range("I122").Value = 338.5
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 3046.5
$ws.Range("L122").Value = 26991
$ws.Range("M122").Value = -596.5
$ws.Range("N122").Value = -31891

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 1190.303
$ws.Range("I131").Value = 1688.8889
$ws.Range("J131").Value = 1003.3333
$ws.Range("K131").Value = 5066.6667
$ws.Range("L131").Value = 3009.9999
$ws.Range("M131").Value = -26.66669999999976
$ws.Range("N131").Value = -13089.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 102949.8
$ws.Range("I102").Value = 1749.8334
$ws.Range("J102").Value = 254749.75
$ws.Range("K102").Value = 1749.8334
$ws.Range("L102").Value = 254749.75
$ws.Range("M102").Value = -127.8334
$ws.Range("N102").Value = -257993.75

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 2699.44
$ws.Range("I126").Value = 1517.1
$ws.Range("K126").Value = 4551.299999999999
$ws.Range("M126").Value = -2081.299999999999

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 4007.6553
$ws.Range("I132").Value = 3770.25
$ws.Range("K132").Value = 11310.75
$ws.Range("M132").Value = -8780.75

# Row 138: Orders Anonymous / White Gold Halfmask of Maiming
$ws.Range("H138").Value = 44083.332
$ws.Range("J138").Value = 44083.332
$ws.Range("L138").Value = 44083.332
$ws.Range("N138").Value = -54363.332

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 10000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 10000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -10224

# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 3085.3333
$ws.Range("I93").Value = 2377.3333
$ws.Range("J93").Value = 7333.3335
$ws.Range("K93").Value = 2377.3333
$ws.Range("L93").Value = 7333.3335
$ws.Range("M93").Value = -1129.3333
$ws.Range("N93").Value = -9829.333500000001

# Row 94: Fitting In / Gaganaskin Hat of Aiming
$ws.Range("H94").Value = 25000
$ws.Range("J94").Value = 25000
$ws.Range("L94").Value = 25000
$ws.Range("N94").Value = -26352

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 30000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -34940

$ws = $wb.Worksheets.Item("WVR")
# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 1461.1765
$ws.Range("I113").Value = 523.75
$ws.Range("J113").Value = 2294.4443
$ws.Range("K113").Value = 1571.25
$ws.Range("L113").Value = 6883.3329
$ws.Range("M113").Value = 598.75
$ws.Range("N113").Value = -11223.3329

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 1252.7931
$ws.Range("I136").Value = 730.0476
$ws.Range("J136").Value = 2625
$ws.Range("K136").Value = 2190.1428
$ws.Range("L136").Value = 7875
$ws.Range("M136").Value = 359.8571999999999
$ws.Range("N136").Value = -12975
